$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "SUREKHA A" / "Teaching Assistant" row (row 9) contents without
# shifting subsequent rows (row 12 stays at row 12).
$ws.Range("A9:B9").ClearContents()

# Update the active selection to A9
$ws.Range("A9").Select()
